$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("D5").Value = 44169
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 5500
$ws.Range("P5").Value = 5750
$ws.Range("S5").Value = 3833

# Row 6 updates
$ws.Range("D6").Value = 44159
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 6500
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 6750
$ws.Range("S6").Value = 4500

# Row 7 updates
$ws.Range("D7").Value = 44176
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 5000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 5500
$ws.Range("S7").Value = 3667
